# Update gh-pages output: increment "想去人数" (F column) counts that changed
# between scrapes for several rows across the 展览 / 演出 / 本地生活 / 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 780
$ws.Range("F3").Value = 2817
$ws.Range("F9").Value = 286
$ws.Range("F11").Value = 11754
$ws.Range("F20").Value = 92
$ws.Range("F23").Value = 3654
$ws.Range("F30").Value = 232
$ws.Range("F32").Value = 311
$ws.Range("F33").Value = 5035
$ws.Range("F35").Value = 1248
$ws.Range("F37").Value = 561
$ws.Range("F39").Value = 546

# --- 演出 (Performances) sheet ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F20").Value = 4

# --- 本地生活 (Local life) sheet ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9067
$ws.Range("F3").Value = 511

# --- 全部类型 (All types) sheet ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 511
$ws.Range("F4").Value = 780
$ws.Range("F5").Value = 2817
$ws.Range("F15").Value = 286
$ws.Range("F17").Value = 11754
$ws.Range("F27").Value = 92
$ws.Range("F30").Value = 3654
$ws.Range("F35").Value = 232
$ws.Range("F40").Value = 1248
$ws.Range("F44").Value = 546
$ws.Range("F45").Value = 4
